$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update Tiempo_Mínimo / Tiempo_Máximo / Tiempo_Promedio values for row 2
$ws.Range("E2").Value = 0.000295941
$ws.Range("F2").Value = 0.022672494
$ws.Range("G2").Value = 0.0005221321025021026

# Update Tiempo_Mínimo / Tiempo_Máximo / Tiempo_Promedio values for row 3
$ws.Range("E3").Value = 0.004604507
$ws.Range("F3").Value = 0.009885973
$ws.Range("G3").Value = 0.005467831646608316
